$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "print preview / multi sheet" help text block in rows 7-11.
# Values are entered in this order so the shared-string table ends up in
# the same order as the target workbook (SUM1/SUM2 etc. already occupy 0-13).
$ws.Range("A9").Value  = "1)Create 1st sheet like this"
$ws.Range("A11").Value = "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas"
$ws.Range("A10").Value = "2)or call macro like in example 070 to save file with a new name (so the entire file would be recreated)"
$ws.Range("A8").Value  = "If you have problems with printing the document:"
$ws.Range("A7").Value  = "Multi sheet issues with print preview"

# Row 11 holds the actual hyperlink (uses the built-in "Hyperlink" style).
# This must be created before the bold heading below so the new font /
# cellXf entries land in the same order as the target workbook.
$ws.Hyperlinks.Add($ws.Range("A11"), "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas")

# Row 7 is a bold heading.
$ws.Range("A7").Font.Bold = $true

# Keep the explicit row height on the new rows (matches the rest of the sheet).
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15

# Page setup for print preview (paper size / orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to G1 like in the final workbook.
[void]$ws.Range("G1").Select()
